$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "**test no BR" test value from G3 (subtitle column for BC02 row)
$ws.Range("G3").ClearContents()

# Update the view/selection state: scroll to show column H as the left-most
# visible column, with M3 as the active/selected cell
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("M3").Select()
